# Correção nos dados: a linha 6 ("grandes regiões e unidades da
# federação") era apenas um cabeçalho de seção sem valores, o que
# deixava todos os dados das regiões/UFs deslocados uma linha para
# baixo (ex.: os valores de "norte" estavam na linha da label
# "rondônia" etc.). Removendo essa linha inteira (com deslocamento das
# linhas abaixo para cima) realinha cada label com os seus valores
# corretos e descarta a última linha, que ficava sem dados.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("6").Delete()
